$d = $word.ActiveDocument

# 1) Strike-through "Battles the chief with nantu and bodiless smile." and
#    "Bodiless smile reveal." -- locate the paragraphs by their text and
#    toggle strike-through on the whole paragraph range (covers both the
#    run and the paragraph mark, matching Word's own behavior when you
#    select a paragraph and press Ctrl+Shift+X / click Strikethrough).
foreach ($par in $d.Paragraphs) {
    $t = $par.Range.Text
    if ($t -eq "Battles the chief with nantu and bodiless smile.`r") {
        $par.Range.Font.StrikeThrough = 1
    }
    elseif ($t -eq "Bodiless smile reveal.`r") {
        $par.Range.Font.StrikeThrough = 1
    }
}

# 2) Insert a new paragraph "Jack sneaks up on Nick during commotion of
#    battle and stabs him with a knife." right before "Frees prisoners
#    during battle", then strike-through that latter paragraph.
foreach ($par in $d.Paragraphs) {
    $t = $par.Range.Text
    if ($t -eq "Frees prisoners during battle`r") {
        $par.Range.InsertParagraphBefore()
        break
    }
}

foreach ($par in $d.Paragraphs) {
    $t = $par.Range.Text
    if ($t -eq "`r") {
        $nxt = $par.Next()
        if ($nxt -ne $null -and $nxt.Range.Text -eq "Frees prisoners during battle`r") {
            $par.Range.Text = "Jack sneaks up on Nick during commotion of battle and stabs him with a knife."
            break
        }
    }
}

foreach ($par in $d.Paragraphs) {
    $t = $par.Range.Text
    if ($t -eq "Frees prisoners during battle`r") {
        $par.Range.Font.StrikeThrough = 1
    }
}

# 3) Remove the stray <w:lastRenderedPageBreak/> rendering hint from the run
#    preceding "Make nick more friendly until the altercation on the
#    ship/planet" by re-typing the run's text (a no-op content change that
#    drops the rendering artifact, same as Word does on edit).
foreach ($par in $d.Paragraphs) {
    $t = $par.Range.Text
    if ($t -eq "Make nick more friendly until the altercation on the ship/planet`r") {
        $par.Range.Text = "Make nick more friendly until the altercation on the ship/planet"
    }
}
